$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": just a selection change (A7:XFD15 -> A7:XFD14)
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Activate()
$wsSummary.Range("A7:XFD14").Select()

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule": add column O (mirrors column N: same format,
# same values) for rows 2 through 15.
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Activate()

$wsRepay.Range("N2:N15").Copy()
$wsRepay.Range("O2:O15").PasteSpecial(-4122)  # xlPasteFormats
$wsRepay.Range("N2:N15").Copy()
$wsRepay.Range("O2:O15").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet "Transactions": update A2:A4 values and change the active
# selection (A2:XFD5 -> D3).
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")
$wsTrans.Activate()

$wsTrans.Range("A2").Value = 115
$wsTrans.Range("A3").Value = 114
$wsTrans.Range("A4").Value = 113

$wsTrans.Range("D3").Select()
